$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Platform Coverage") ---
# Row 5 (Treatment/Campaign/MDA, 50-65) was a duplicate of the age-band
# series and gets removed; rows 6-9 shift up to become rows 5-8.
$ws1.Rows.Item(5).Delete()

# Row 4 (now the last MDA row) absorbs the widened age range.
$ws1.Range("G4").Value = 65

# New empty, center-aligned cell next to it.
$ws1.Range("H4").HorizontalAlignment = -4108  # xlCenter

# --- Sheet2 ("MarketShare") view/selection update ---
$ws2.Activate()
[void]$ws2.Range("N22").Select()

# Restore Sheet1 as the active tab, with its new selection.
$ws1.Activate()
[void]$ws1.Range("H4").Select()
